$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New city rows to append (rows 185-188), mirroring the existing layout:
# A = City English Name, B = City English Name (repeat), C = City Arabic Name,
# D = Latitude, E = Longitude, F = Area (Arabic), G = Region (Arabic)
$qassim = "منطقة" + [char]0x00A0 + "القصيم"

$newRows = @(
    @("Ash Shinan",      "الشنان",     27.176233, 42.443098999999997, "منطقة حائل",   "شمال المملكة"),
    @("Badr Al Janoub",  "بدر الجنوب", 17.879346999999999, 43.719791999999998, "منطقة نجران",  "جنوب المملكة"),
    @("An Nabhaniyah",   "النبهانية",  25.85782, 43.067594999999997, $qassim, "وسط المملكة"),
    @("Dariyah",         "ضرية",       24.722766, 42.932467000000003, $qassim, "وسط المملكة")
)

$startRow = 185
$templateRow = 184
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Carry the existing data-row formatting (thin border style) onto the new row.
    $src = $ws.Range("A" + $templateRow + ":G" + $templateRow)
    $dst = $ws.Range("A" + $r + ":G" + $r)
    $src.Copy($dst)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

# Extend the sheet's stored selection to cover the new data extent (A1:G188),
# matching what Excel records after the rows were appended.
[void]$ws.Range("A1:G188").Select()
